$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "10/19/2025"
$ws.Range("A63").ClearFormats()
$ws.Range("B63").Value = 9486.57
